$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 263.66666
$ws.Range("I5").Value = 263.66666
$ws.Range("K5").Value = 263.66666
$ws.Range("M5").Value = -148.66666

$ws.Range("H11").Value = 166672000
$ws.Range("I11").Value = 166672000
$ws.Range("K11").Value = 166672000
$ws.Range("M11").Value = -166671860

$ws.Range("H16").Value = 3751.125
$ws.Range("J16").Value = 4001.5
$ws.Range("L16").Value = 4001.5
$ws.Range("N16").Value = -4461.5

$ws.Range("H86").Value = 10210.211
$ws.Range("I86").Value = 10352.647
$ws.Range("K86").Value = 10352.647
$ws.Range("M86").Value = -9229.647000000001

$ws.Range("H89").Value = 10210.211
$ws.Range("I89").Value = 10352.647
$ws.Range("K89").Value = 51763.235
$ws.Range("M89").Value = -46147.235

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1025246.25
$ws.Range("I31").Value = 27401.5
$ws.Range("K31").Value = 27401.5
$ws.Range("M31").Value = -27106.5

$ws.Range("H34").Value = 1025246.25
$ws.Range("I34").Value = 27401.5
$ws.Range("K34").Value = 27401.5
$ws.Range("M34").Value = -27199.5

$ws.Range("H74").Value = 17400
$ws.Range("J74").Value = 17400
$ws.Range("L74").Value = 17400
$ws.Range("N74").Value = -19148

$ws.Range("H77").Value = 17400
$ws.Range("J77").Value = 17400
$ws.Range("L77").Value = 52200
$ws.Range("N77").Value = -60936

$ws.Range("H86").Value = 3000
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1877
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 3000
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -9384
$ws.Range("N89").ClearContents()

$ws.Range("H93").Value = 62499.5
$ws.Range("I93").Value = 16999
$ws.Range("J93").Value = 108000
$ws.Range("K93").Value = 16999
$ws.Range("L93").Value = 108000
$ws.Range("M93").Value = -15127
$ws.Range("N93").Value = -111744

$ws.Range("H132").Value = 312
$ws.Range("I132").Value = 312
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 936
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 1594
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 1007680.9
$ws.Range("I134").Value = 2501475
$ws.Range("K134").Value = 7504425
$ws.Range("M134").Value = -7501890

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 5000.8
$ws.Range("I133").Value = 4778.6665
$ws.Range("K133").Value = 14335.9995
$ws.Range("M133").Value = -9275.999500000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 4098.8335
$ws.Range("I3").Value = 20000
$ws.Range("J3").Value = 918.6
$ws.Range("K3").Value = 20000
$ws.Range("L3").Value = 918.6
$ws.Range("M3").Value = -19884
$ws.Range("N3").Value = -1150.6

$ws.Range("H5").Value = 11725.454
$ws.Range("I5").Value = 8886.111000000001
$ws.Range("J5").Value = 24502.5
$ws.Range("K5").Value = 8886.111000000001
$ws.Range("L5").Value = 24502.5
$ws.Range("M5").Value = -8774.111000000001
$ws.Range("N5").Value = -24726.5

$ws.Range("H10").Value = 3966.6667
$ws.Range("I10").Value = 3450
$ws.Range("J10").Value = 5000
$ws.Range("K10").Value = 3450
$ws.Range("L10").Value = 5000
$ws.Range("M10").Value = -3281
$ws.Range("N10").Value = -5338

$ws.Range("H11").Value = 8145851.5
$ws.Range("I11").Value = 9726540
$ws.Range("J11").Value = 6193235.5
$ws.Range("K11").Value = 9726540
$ws.Range("L11").Value = 6193235.5
$ws.Range("M11").Value = -9726401
$ws.Range("N11").Value = -6193513.5

$ws.Range("H101").Value = 69999.5
$ws.Range("J101").Value = 69999.5
$ws.Range("L101").Value = 69999.5
$ws.Range("N101").Value = -76489.5

$ws.Range("H102").Value = 3530.8918
$ws.Range("I102").Value = 2884.4482
$ws.Range("K102").Value = 2884.4482
$ws.Range("M102").Value = -1262.4482

$ws.Range("H141").Value = 25000
$ws.Range("J141").Value = 25000
$ws.Range("L141").Value = 25000
$ws.Range("N141").Value = -35360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 38702.586
$ws.Range("I7").Value = 3792.7144
$ws.Range("J7").Value = 130341
$ws.Range("K7").Value = 3792.7144
$ws.Range("L7").Value = 130341
$ws.Range("M7").Value = -3680.7144
$ws.Range("N7").Value = -130565

$ws.Range("H16").Value = 2150
$ws.Range("I16").Value = 1709.3334
$ws.Range("J16").Value = 3031.3333
$ws.Range("K16").Value = 1709.3334
$ws.Range("L16").Value = 3031.3333
$ws.Range("M16").Value = -1539.3334
$ws.Range("N16").Value = -3371.3333

$ws.Range("H22").Value = 2460.7407
$ws.Range("I22").Value = 2472
$ws.Range("K22").Value = 2472
$ws.Range("M22").Value = -2177

$ws.Range("H27").Value = 2460.7407
$ws.Range("I27").Value = 2472
$ws.Range("K27").Value = 2472
$ws.Range("M27").Value = -2365

$ws.Range("H46").Value = 3871.4075
$ws.Range("I46").Value = 3298.8333
$ws.Range("K46").Value = 3298.8333
$ws.Range("M46").Value = -3110.8333

$ws.Range("H61").Value = 1001
$ws.Range("I61").Value = 1001
$ws.Range("K61").Value = 1001
$ws.Range("M61").Value = -799

$ws.Range("H103").Value = 35934
$ws.Range("J103").Value = 35934
$ws.Range("L103").Value = 35934
$ws.Range("N103").Value = -38278

$ws.Range("H113").Value = 1001
$ws.Range("I113").Value = 1001
$ws.Range("K113").Value = 1001
$ws.Range("M113").Value = 1169

$ws.Range("H126").Value = 38702.586
$ws.Range("I126").Value = 3792.7144
$ws.Range("J126").Value = 130341
$ws.Range("K126").Value = 11378.1432
$ws.Range("L126").Value = 391023
$ws.Range("M126").Value = -8908.143199999999
$ws.Range("N126").Value = -395963

$ws.Range("H127").Value = 95630
$ws.Range("J127").Value = 95630
$ws.Range("L127").Value = 95630
$ws.Range("N127").Value = -105550

$ws.Range("H132").Value = 68717.03
$ws.Range("I132").Value = 39824.58
$ws.Range("K132").Value = 119473.74
$ws.Range("M132").Value = -116943.74

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H102").Value = 110000
$ws.Range("J102").Value = 110000
$ws.Range("L102").Value = 110000
$ws.Range("N102").Value = -116490

$ws.Range("H132").Value = 8556.071
$ws.Range("I132").Value = 1525.7273
$ws.Range("J132").Value = 34334
$ws.Range("K132").Value = 4577.1819
$ws.Range("L132").Value = 103002
$ws.Range("M132").Value = -2047.1819
$ws.Range("N132").Value = -108062

$ws.Range("H139").Value = 49000
$ws.Range("J139").Value = 49000
$ws.Range("L139").Value = 49000
$ws.Range("N139").Value = -59280
